$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The elective-course table was missing a row for course "GS2515"
# (which has two offered course-name variants). Insert two new rows
# at row 15 for GS2515, then drop the now-duplicated GS2615 rows that
# got shifted down to rows 32:33 (their course-name text was reused
# for the new GS2515 rows).

$ws.Rows("15:16").Insert()

$ws.Range("A15").Value = "GS2515"
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 2021
$ws.Range("D15").Value = "사회주의: 이론과 역사"

$ws.Range("A16").Value = "GS2515"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 2021
$ws.Range("D16").Value = "자본주의와 사회주의의 역사"

$ws.Rows("32:33").Delete()

# Restore the saved view/selection state as recorded after this edit.
$ws.Range("I20").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
